$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 248, shifting rows 248:343 down to 249:344
$ws.Rows.Item(248).Insert()

# Populate the new row 248 with data
$ws.Cells.Item(248, 1).Value = 4
$ws.Cells.Item(248, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(248, 3).Value = "Los Lagos"
$ws.Cells.Item(248, 4).Value = 44784
$ws.Cells.Item(248, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(248, 5).Value = 10
$ws.Cells.Item(248, 6).Value = 100112045
$ws.Cells.Item(248, 7).Value = "Zapallo"
$ws.Cells.Item(248, 8).Value = "Paine"
$ws.Cells.Item(248, 9).Value = "1a (guarda)"
$ws.Cells.Item(248, 10).Value = 500
$ws.Cells.Item(248, 11).Value = 550
$ws.Cells.Item(248, 12).Value = 650
$ws.Cells.Item(248, 13).Value = 600
$ws.Cells.Item(248, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(248, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(248, 16).Value = 600
$ws.Cells.Item(248, 17).Value = 1
$ws.Cells.Item(248, 18).Value = "Hortaliza"
